{"js": "// The diff merges several runs inside two list-paragraphs into a single\n// run each. The visible text does not change - only the run boundaries\n// go away. We find the two paragraphs by matching their current\n// (pre-edit) full text against the two known sentences, then rewrite\n// each paragraph's content using its own already-loaded text (so any\n// special characters, e.g. a non-breaking space, survive untouched) via\n// insertText(..., \"Replace\"), which merges the paragraph down to one run\n// while keeping the paragraph/list formatting intact.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Vyjad\u0159uje opr\u00e1vn\u011bn\u00fd v\u00edcen\u00e1klad dopravce\",\n  \"PAX nem\u00e1 zakoupenou j\u00edzdenku po n\u00e1stupu do vlaku v obsazen\u00e9 stanici a dopravce nestanovil, \u017ee se p\u0159ir\u00e1\u017eka neplat\u00ed\",\n];\n\n// Collapse any run of whitespace (regular space, non-breaking space, ...)\n// down to a single space so matching does not depend on exactly which\n// whitespace character the source document happens to use.\nconst normalizeWhitespace = (text) => text.replace(/\\s+/g, \" \");\n\nconst targetParagraphs = paragraphs.items.filter((paragraph) =>\n  targetTexts.some(\n    (text) => normalizeWhitespace(paragraph.text) === normalizeWhitespace(text)\n  )\n);\n\nif (targetParagraphs.length !== targetTexts.length) {\n  throw new Error(\n    \"Expected to find \" +\n      targetTexts.length +\n      \" target paragraphs, found \" +\n      targetParagraphs.length\n  );\n}\n\nfor (const paragraph of targetParagraphs) {\n  // Re-use the exact text already present in the paragraph (the\n  // concatenation of its current runs) so the merge is byte-for-byte\n  // faithful, then collapse it into a single run.\n  paragraph.insertText(paragraph.text, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The diff merges several runs inside two list-paragraphs into a single\n# run each; the visible text itself does not change. We locate the two\n# paragraphs by their current text (normalizing whitespace so it does not\n# matter whether the source uses a regular or a non-breaking space), then\n# rewrite each paragraph's range using the text the paragraph already\n# contains. Assigning that text (without the trailing paragraph mark) back\n# to a fresh Range spanning the whole paragraph (mark included) collapses\n# all of the paragraph's runs into a single run while Word automatically\n# keeps the paragraph mark, style and numbering untouched.\n\n$d = $word.ActiveDocument\n\n$targets = @(\n    \"Vyjad\u0159uje opr\u00e1vn\u011bn\u00fd v\u00edcen\u00e1klad dopravce\",\n    \"PAX nem\u00e1 zakoupenou j\u00edzdenku po n\u00e1stupu do vlaku v obsazen\u00e9 stanici a dopravce nestanovil, \u017ee se p\u0159ir\u00e1\u017eka neplat\u00ed\"\n)\n\nfunction Normalize-Whitespace($text) {\n    return ($text -replace '\\s+', ' ')\n}\n\n$matchCount = 0\n\nforeach ($p in $d.Paragraphs) {\n    $paraRange = $p.Range\n    $fullText = $paraRange.Text\n    # Strip the trailing paragraph mark (CR) / cell mark before comparing.\n    $text = $fullText.TrimEnd([char]13, [char]7)\n    $normalized = Normalize-Whitespace $text\n\n    foreach ($target in $targets) {\n        if ($normalized -eq (Normalize-Whitespace $target)) {\n            $matchCount = $matchCount + 1\n            $editRange = $d.Range($paraRange.Start, $paraRange.End)\n            $editRange.Text = $text\n            break\n        }\n    }\n}\n\nif ($matchCount -ne $targets.Count) {\n    throw \"Expected to find $($targets.Count) target paragraphs, found $matchCount\"\n}\n"}
